# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# figures for each coin row on the sheet, matching the latest scrape.
#
# Note on column D: most "Price" strings look like plain decimals
# (e.g. "21.20", "7.500", "0.02380") and Excel's normal type inference
# would silently reinterpret them as numbers -- dropping the significant
# trailing zeros we need to keep literally. To avoid that, those cells
# are briefly switched to a Text number format before the value is
# written, then restored to the workbook's default ("Normal") style so
# no stray formatting is left behind. Multi-dot price strings (e.g.
# "28.508.77") can never be parsed as a number, so they're set directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.508.77"

$ws.Range("D3").Value = "1.828.01"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5161"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08319"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.00%  "

$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.86%  "

$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.500"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "1.826.34"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001121"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("D23").Value = "28.554.32"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("D28").Value = "2.035.20"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.100"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.743"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.63%  "

$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2235"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02380"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.291"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.774"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.38%  "

$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.395"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6167"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.806"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.208"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06985"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.35"
$ws.Range("D51").Style = "Normal"
